# Weir Calibration Field Form 2019 - "starting on data summary"
#
# This script:
#   1. Updates the row-visibility (hidden) state of a set of response rows
#      on the "Form Responses 1" sheet (spot-checking / cleaning rows as
#      part of starting the data summary).
#   2. Corrects a mis-keyed data value in H252 (50 -> 5).
#   3. Re-applies AutoFilter: clears the old "SITE ID = SDR-751" filter and
#      instead filters column H ("Height above (or below) v-notch (cm)")
#      down to a specific set of values, over the full used range A1:AA673.
#   4. Updates the sheet's scroll position / selected cell to reflect where
#      the user was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Form Responses 1")

# --- 1. Rebuild the AutoFilter -------------------------------------------

# Drop the previous filter entirely before re-applying a new one so the
# old "SITE ID" criteria doesn't linger alongside the new column. Do this
# first (while H252 is still its original 50) so the new filter's
# row-visibility pass matches what was recorded when the filter was
# actually applied, before the data-entry fix below.
$ws.AutoFilterMode = $false

$ws.Range("A1:AA673").AutoFilter(8, @("14.5", "15", "16.8", "17.5", "21.5", "50"), 7)

# Keep the named filter-database range in sync with the new AutoFilter range.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Form Responses 1'!`$A`$1:`$AA`$673"
    }
}

# --- 2. Spot-fix a handful of rows' visibility ---------------------------

$rowsToHide = @(16, 80, 133, 136, 221, 243, 322, 385, 424, 469, 505, 533, 554, 618, 643, 645)
foreach ($r in $rowsToHide) {
    $ws.Rows.Item($r).Hidden = $true
}

$rowsToShow = @(36, 148, 219, 224, 252, 497)
foreach ($r in $rowsToShow) {
    $ws.Rows.Item($r).Hidden = $false
}

# --- 3. Fix a data-entry value --------------------------------------------

$ws.Range("H252").Value = 5

# --- 4. Update the view's frozen-pane scroll position / selection -------

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("H497").Select()
